$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.806.67'
$ws.Range('E2').Value = '  +4.72%  '
$ws.Range('D3').Value = '2.281.40'
$ws.Range('E3').Value = '  +2.22%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '''231.18'
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('D6').Value = '''0.628'
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('D7').Value = '''64.01'
$ws.Range('E7').Value = '  +5.26%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '''0.423'
$ws.Range('E9').Value = '  +4.02%  '
$ws.Range('D10').Value = '''0.0949'
$ws.Range('E10').Value = '  +4.71%  '
$ws.Range('D11').Value = '''57.70'
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('D12').Value = '''26.42'
$ws.Range('E12').Value = '  +16.35%  '
$ws.Range('E13').Value = '  +0.09%  '
$ws.Range('D14').Value = '2.622.69'
$ws.Range('E14').Value = '  +2.25%  '
$ws.Range('D15').Value = '''15.72'
$ws.Range('E15').Value = '  -0.17%  '
$ws.Range('D16').Value = '''5.90'
$ws.Range('E16').Value = '  +4.80%  '
$ws.Range('D17').Value = '''0.812'
$ws.Range('E17').Value = '  +0.84%  '
$ws.Range('D18').Value = '2.285.65'
$ws.Range('E18').Value = '  +2.00%  '
$ws.Range('D19').Value = '43.736.15'
$ws.Range('E19').Value = '  +4.64%  '
$ws.Range('D20').Value = '0.0₃0943'
$ws.Range('E20').Value = '  +4.01%  '
$ws.Range('D21').Value = '''73.18'
$ws.Range('E21').Value = '  +0.81%  '
$ws.Range('D22').Value = '''6.16'
$ws.Range('E22').Value = '  +0.47%  '
$ws.Range('D23').Value = '''249.92'
$ws.Range('E23').Value = '  +0.69%  '
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').Value = '''2.56'
$ws.Range('E25').Value = '  +7.46%  '
$ws.Range('D26').Value = '''2.33'
$ws.Range('E26').Value = '  +0.95%  '
$ws.Range('D27').Value = '''9.83'
$ws.Range('E27').Value = '  +1.28%  '
$ws.Range('D28').Value = '''171.47'
$ws.Range('D29').Value = '''0.138'
$ws.Range('E29').Value = '  -3.77%  '
$ws.Range('D30').Value = '''20.50'
$ws.Range('E30').Value = '  +2.66%  '
$ws.Range('D31').Value = '''1.43'
$ws.Range('E31').Value = '  +1.79%  '
$ws.Range('D32').Value = '''2.75'
$ws.Range('E32').Value = '  +3.15%  '
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('D34').Value = '''0.0698'
$ws.Range('E34').Value = '  +6.54%  '
$ws.Range('D35').Value = '''5.12'
$ws.Range('E35').Value = '  +0.88%  '
$ws.Range('D36').Value = '''4.71'
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('D37').Value = '''6.61'
$ws.Range('E37').Value = '  +0.44%  '
$ws.Range('D38').Value = '''3.70'
$ws.Range('E38').Value = '  +1.54%  '
$ws.Range('D39').Value = '''2.36'
$ws.Range('E39').Value = '  -1.74%  '
$ws.Range('D40').Value = '''0.0248'
$ws.Range('E40').Value = '  +2.94%  '
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('E42').Value = '  +8.97%  '
$ws.Range('D43').Value = '''10.92'
$ws.Range('E43').Value = '  +26.57%  '
$ws.Range('D44').Value = '''8.49'
$ws.Range('E44').Value = '  -1.55%  '
$ws.Range('D45').Value = '''0.000220'
$ws.Range('E45').Value = '  -8.52%  '
$ws.Range('D46').Value = '''1.21'
$ws.Range('E46').Value = '  -1.24%  '
$ws.Range('D47').Value = '''0.0963'
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('D48').Value = '''97.73'
$ws.Range('E48').Value = '  -1.08%  '
$ws.Range('D49').Value = '1.482.94'
$ws.Range('E49').Value = '  +0.68%  '
$ws.Range('D50').Value = '''16.83'
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('D51').Value = '''2.35'
$ws.Range('E51').Value = '  +3.45%  '
